$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "259.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.61%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.19%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.674"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.33%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06029"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.40%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.667"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8589"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.14%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9319"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.57%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.77%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04614"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "20.43%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07007"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.09%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03116"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.77%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09135"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.09%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001532"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.24%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006068"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.05%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006133"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.42%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.465"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.155"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.51%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.74%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.24%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1295"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.28%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.134"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.36%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04233"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.50%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001214"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.59%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004050"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.85%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001197"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.16%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "13.62%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03848"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.05%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1117"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.25%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003898"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-37.67%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002414"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.57%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "29.28%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.47%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.15%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-16.67%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1303"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.37%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002095"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.15%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001995"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.15%"
